$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 "29.815.31"
$ws.Range("E2").Value = "  -2.65%  "
Set-TextValue 3 4 "2.094.74"
$ws.Range("E3").Value = "  -0.85%  "
Set-TextValue 4 4 "1.012"
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue 5 4 "344.50"
$ws.Range("E5").Value = "  +1.73%  "
Set-TextValue 6 4 "1.011"
$ws.Range("E6").Value = "  -0.16%  "
Set-TextValue 7 4 "0.5174"
$ws.Range("E7").Value = "  -1.46%  "
Set-TextValue 8 4 "0.4464"
$ws.Range("E8").Value = "  -0.96%  "
Set-TextValue 9 4 "0.09506"
$ws.Range("E9").Value = "  +5.29%  "
Set-TextValue 10 4 "52.02"
$ws.Range("E10").Value = "  -2.83%  "
Set-TextValue 11 4 "1.170"
$ws.Range("E11").Value = "  +0.15%  "
Set-TextValue 12 4 "25.05"
$ws.Range("E12").Value = "  +2.84%  "
Set-TextValue 13 4 "2.105.83"
$ws.Range("E13").Value = "  -0.94%  "
Set-TextValue 14 4 "6.722"
$ws.Range("E14").Value = "  -0.86%  "
Set-TextValue 15 4 "8.044"
$ws.Range("E15").Value = "  -0.19%  "
Set-TextValue 16 4 "99.30"
$ws.Range("E16").Value = "  +1.44%  "
Set-TextValue 17 4 "0.00001163"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("E18").Value = "  -0.35%  "
Set-TextValue 19 4 "0.06708"
$ws.Range("E19").Value = "  +0.07%  "
Set-TextValue 20 4 "20.51"
Set-TextValue 21 4 "1.009"
$ws.Range("E21").Value = "  -0.29%  "
Set-TextValue 22 4 "6.174"
$ws.Range("E22").Value = "  -2.35%  "
Set-TextValue 23 4 "29.981.23"
$ws.Range("E23").Value = "  -2.43%  "
Set-TextValue 24 4 "12.64"
$ws.Range("E24").Value = "  -1.21%  "
Set-TextValue 25 4 "2.324"
$ws.Range("E25").Value = "  -2.17%  "
Set-TextValue 26 4 "2.356.38"
$ws.Range("E26").Value = "  -0.73%  "
Set-TextValue 27 4 "22.00"
$ws.Range("E27").Value = "  -1.50%  "
Set-TextValue 28 4 "163.58"
$ws.Range("E28").Value = "  -1.07%  "
Set-TextValue 29 4 "2.526"
$ws.Range("E29").Value = "  -0.42%  "
Set-TextValue 30 4 "133.59"
$ws.Range("E30").Value = "  -0.92%  "
Set-TextValue 31 4 "1.151"
$ws.Range("E31").Value = "  -3.56%  "
Set-TextValue 32 4 "0.1055"
$ws.Range("E32").Value = "  -1.70%  "
Set-TextValue 33 4 "1.610"
$ws.Range("E33").Value = "  -1.35%  "
Set-TextValue 34 4 "6.210"
$ws.Range("E34").Value = "  -2.46%  "
Set-TextValue 35 4 "3.954"
$ws.Range("E35").Value = "  +0.17%  "
Set-TextValue 36 4 "6.139"
$ws.Range("E36").Value = "  +4.33%  "
Set-TextValue 37 4 "10.08"
$ws.Range("E37").Value = "  -2.02%  "
Set-TextValue 38 4 "0.02565"
$ws.Range("E38").Value = "  -3.22%  "
Set-TextValue 39 4 "0.06733"
$ws.Range("E39").Value = "  -1.42%  "
Set-TextValue 40 4 "0.2272"
$ws.Range("E40").Value = "  -1.92%  "
Set-TextValue 43 4 "1.309"
$ws.Range("E43").Value = "  +3.80%  "
Set-TextValue 44 4 "0.6643"
$ws.Range("E44").Value = "  +3.44%  "
Set-TextValue 45 4 "14.18"
$ws.Range("E45").Value = "  -5.57%  "
Set-TextValue 46 4 "2.272"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("E47").Value = "  -1.69%  "
Set-TextValue 48 4 "1.217"
$ws.Range("E48").Value = "  -2.94%  "
Set-TextValue 49 4 "0.00000000338"
$ws.Range("E49").Value = "  -8.15%  "
Set-TextValue 50 4 "81.49"
$ws.Range("E50").Value = "  -1.60%  "
Set-TextValue 51 4 "0.07143"
$ws.Range("E51").Value = "  -2.23%  "

# Row 41/42 swap (TheSandbox <-> Aptos)
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue 41 4 "12.42"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue 42 4 "0.6873"
$ws.Range("E42").Value = "  +0.05%  "
